$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the rows that no longer exist (rows 4-7 from the old 6x7 data block) ---
$ws.Range("A4:F7").ClearContents()

# --- Row 1 (headers) ---
# Write "same" before the distance/out headers so the shared-string table append
# order matches the target layout exactly.
$ws.Range("A1").Value = "down"
$ws.Range("B1").Value = "right"
$ws.Range("C1").Value = "up"
$ws.Range("D1").Value = "left"
$ws.Range("K1").Value = "same"
$ws.Range("E1").Value = "distance to Food"
$ws.Range("F1").Value = "distance to Wall"
$ws.Range("H1").Value = "out"
$ws.Range("J1").Value = "left"
$ws.Range("L1").Value = "right"

# --- Row 2 ---
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 20
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0

# --- Row 3 ---
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 18
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0

# --- Column widths for D:E ---
$ws.Range("D1:E1").ColumnWidth = 13.498697916666666

# --- Selection ---
[void]$ws.Range("K6").Select()

# --- Window position ---
$excel.Windows.Item(1).Left = 4800
$excel.Windows.Item(1).Top = 2620
